# Replace the four "Waktu Kampanye ... 2018 untuk Perseus: ..." paragraphs
# with a single plain run containing the translated Cygnus campaign dates.

$d = $word.ActiveDocument

$oldMarker = "Waktu Kampanye"
$oldTail   = "2018 untuk Perseus: 30 Oktober-8 November dan 29 November-8 Desember"
$newText   = "Waktu Kampanye Cygnus: 10-19 Agustus, 9-18 September, 8-17 Oktober"

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs($i)
    $full = $p.Range.Text

    if ($full.Contains($oldMarker) -and $full.Contains($oldTail)) {
        $r = $p.Range

        # Exclude the trailing paragraph mark from the range we touch.
        $contentEnd = $r.End - 1
        $target = $d.Range($r.Start, $contentEnd)

        # Wipe out every run (and their formatting) in the paragraph...
        $target.Delete()

        # ...then insert a single plain run with the new combined text.
        $insertionPoint = $d.Range($r.Start, $r.Start)
        $insertionPoint.InsertBefore($newText)
    }
}
